$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.471.32"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "2.108.74"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5271"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09017"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.26%  "

$ws.Range("D13").Value = "2.098.35"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.784"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.816"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001130"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06622"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.316"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").Value = "30.522.66"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("D26").Value = "2.347.00"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.585"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.198"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.668"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.163"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.924"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06836"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.580"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2303"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6919"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.75%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6397"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.663"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.34%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.249"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
